# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the f6dc4ae2-... file row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 is the f6dc4ae2-... file; column G is "Latest HO Xliff Generate Date"
$overview.Range("G3").Value = "2016-08-27 04:46:47"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 is the f6dc4ae2-... file.
# Column H = "Correspond Handoff Datetime"
$zhcn.Range("H3").Value = "2016-08-27 04:46:43"
# Column K = "Correspond Handback DateTime"
$zhcn.Range("K3").Value = "2016-08-27 04:46:59"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 is the f6dc4ae2-... file.
# Column H = "Correspond Handoff Datetime"
$dede.Range("H3").Value = "2016-08-27 04:46:47"
# Column K = "Correspond Handback DateTime"
$dede.Range("K3").Value = "2016-08-27 04:47:10"
